# Appends 4 new fixtures (rows 39-42) to the croatia/hnl/2023-2024 sheet.
# New rows copy the formatting (styles) of the final existing data row,
# then their values are overwritten with the new match data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

# --- New fixture row 1 ---
$srcRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, $lastCol))
$newRow = $lastRow + 1
$dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, $lastCol))
$srcRange.Copy($dstRange)

$ws.Cells.Item($newRow, 1).Value = 38
$ws.Cells.Item($newRow, 2).Value = "croatia"
$ws.Cells.Item($newRow, 3).Value = "hnl"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45192.70833333334
$ws.Cells.Item($newRow, 6).Value = "Istra 1961"
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = "Gorica"
$ws.Cells.Item($newRow, 9).Value = 1
$ws.Cells.Item($newRow, 10).Value = 2.21
$ws.Cells.Item($newRow, 11).Value = "17/09/2023 18:43"
$ws.Cells.Item($newRow, 12).Value = 2.71
$ws.Cells.Item($newRow, 13).Value = "23/09/2023 16:57"
$ws.Cells.Item($newRow, 14).Value = 3.35
$ws.Cells.Item($newRow, 15).Value = "17/09/2023 18:43"
$ws.Cells.Item($newRow, 16).Value = 3.06
$ws.Cells.Item($newRow, 17).Value = "23/09/2023 16:57"
$ws.Cells.Item($newRow, 18).Value = 3.12
$ws.Cells.Item($newRow, 19).Value = "17/09/2023 18:43"
$ws.Cells.Item($newRow, 20).Value = 2.88
$ws.Cells.Item($newRow, 21).Value = "23/09/2023 16:57"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/croatia/hnl/istra-1961-hnk-gorica/CC5UERDE/"

$lastRow = $newRow

# --- New fixture row 2 ---
$srcRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, $lastCol))
$newRow = $lastRow + 1
$dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, $lastCol))
$srcRange.Copy($dstRange)

$ws.Cells.Item($newRow, 1).Value = 39
$ws.Cells.Item($newRow, 2).Value = "croatia"
$ws.Cells.Item($newRow, 3).Value = "hnl"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45192.80208333334
$ws.Cells.Item($newRow, 6).Value = "Hajduk Split"
$ws.Cells.Item($newRow, 7).Value = 1
$ws.Cells.Item($newRow, 8).Value = "Lok. Zagreb"
$ws.Cells.Item($newRow, 9).Value = 0
$ws.Cells.Item($newRow, 10).Value = 1.34
$ws.Cells.Item($newRow, 11).Value = "16/09/2023 20:42"
$ws.Cells.Item($newRow, 12).Value = 1.34
$ws.Cells.Item($newRow, 13).Value = "23/09/2023 19:06"
$ws.Cells.Item($newRow, 14).Value = 5.13
$ws.Cells.Item($newRow, 15).Value = "16/09/2023 20:42"
$ws.Cells.Item($newRow, 16).Value = 5.04
$ws.Cells.Item($newRow, 17).Value = "23/09/2023 19:13"
$ws.Cells.Item($newRow, 18).Value = 8.24
$ws.Cells.Item($newRow, 19).Value = "16/09/2023 20:42"
$ws.Cells.Item($newRow, 20).Value = 9.35
$ws.Cells.Item($newRow, 21).Value = "23/09/2023 19:13"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/croatia/hnl/hajduk-split-lok-zagreb/b91QF7b8/"

$lastRow = $newRow

# --- New fixture row 3 ---
$srcRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, $lastCol))
$newRow = $lastRow + 1
$dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, $lastCol))
$srcRange.Copy($dstRange)

$ws.Cells.Item($newRow, 1).Value = 40
$ws.Cells.Item($newRow, 2).Value = "croatia"
$ws.Cells.Item($newRow, 3).Value = "hnl"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45193.70833333334
$ws.Cells.Item($newRow, 6).Value = "Slaven Belupo"
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = "Rijeka"
$ws.Cells.Item($newRow, 9).Value = 1
$ws.Cells.Item($newRow, 10).Value = 4.47
$ws.Cells.Item($newRow, 11).Value = "17/09/2023 20:42"
$ws.Cells.Item($newRow, 12).Value = 4.8
$ws.Cells.Item($newRow, 13).Value = "24/09/2023 16:59"
$ws.Cells.Item($newRow, 14).Value = 3.66
$ws.Cells.Item($newRow, 15).Value = "17/09/2023 20:42"
$ws.Cells.Item($newRow, 16).Value = 4.03
$ws.Cells.Item($newRow, 17).Value = "24/09/2023 16:59"
$ws.Cells.Item($newRow, 18).Value = 1.72
$ws.Cells.Item($newRow, 19).Value = "17/09/2023 20:42"
$ws.Cells.Item($newRow, 20).Value = 1.69
$ws.Cells.Item($newRow, 21).Value = "24/09/2023 16:36"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/croatia/hnl/slaven-belupo-rijeka/6FEvh461/"

$lastRow = $newRow

# --- New fixture row 4 ---
$srcRange = $ws.Range($ws.Cells.Item($lastRow, 1), $ws.Cells.Item($lastRow, $lastCol))
$newRow = $lastRow + 1
$dstRange = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, $lastCol))
$srcRange.Copy($dstRange)

$ws.Cells.Item($newRow, 1).Value = 41
$ws.Cells.Item($newRow, 2).Value = "croatia"
$ws.Cells.Item($newRow, 3).Value = "hnl"
$ws.Cells.Item($newRow, 4).Value = "2023-2024"
$ws.Cells.Item($newRow, 5).Value = 45193.80208333334
$ws.Cells.Item($newRow, 6).Value = "Rudes"
$ws.Cells.Item($newRow, 7).Value = 1
$ws.Cells.Item($newRow, 8).Value = "D. Zagreb"
$ws.Cells.Item($newRow, 9).Value = 5
$ws.Cells.Item($newRow, 10).Value = 10.1
$ws.Cells.Item($newRow, 11).Value = "20/09/2023 02:12"
$ws.Cells.Item($newRow, 12).Value = 12.27
$ws.Cells.Item($newRow, 13).Value = "24/09/2023 19:13"
$ws.Cells.Item($newRow, 14).Value = 5.54
$ws.Cells.Item($newRow, 15).Value = "20/09/2023 02:12"
$ws.Cells.Item($newRow, 16).Value = 6.23
$ws.Cells.Item($newRow, 17).Value = "24/09/2023 19:13"
$ws.Cells.Item($newRow, 18).Value = 1.28
$ws.Cells.Item($newRow, 19).Value = "20/09/2023 02:12"
$ws.Cells.Item($newRow, 20).Value = 1.23
$ws.Cells.Item($newRow, 21).Value = "24/09/2023 19:13"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/croatia/hnl/rudes-din-zagreb/p6DriOL7/"

$lastRow = $newRow

